$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns retain text formatting so
# numeric-looking strings (e.g. "6.64") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '59.750.79'
$ws.Range("E2").Value = '  +2.67%  '
$ws.Range("D3").Value = '2.600.28'
$ws.Range("E3").Value = '  +1.17%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '567.83'
$ws.Range("E5").Value = '  +5.04%  '
$ws.Range("D6").Value = '143.60'
$ws.Range("E6").Value = '  +1.20%  '
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = '0.599'
$ws.Range("E8").Value = '  +3.94%  '
$ws.Range("D9").Value = '2.610.35'
$ws.Range("E9").Value = '  +1.66%  '
$ws.Range("D10").Value = '6.64'
$ws.Range("E10").Value = '  -1.75%  '
$ws.Range("D11").Value = '0.103'
$ws.Range("E11").Value = '  +3.71%  '
$ws.Range("D12").Value = '0.152'
$ws.Range("E12").Value = '  +10.52%  '
$ws.Range("D13").Value = '0.341'
$ws.Range("E13").Value = '  +3.32%  '
$ws.Range("D14").Value = '3.064.92'
$ws.Range("E14").Value = '  +1.28%  '
$ws.Range("D15").Value = '59.815.23'
$ws.Range("E15").Value = '  +2.92%  '
$ws.Range("D16").Value = '21.94'
$ws.Range("E16").Value = '  +7.27%  '
$ws.Range("D17").Value = '0.0000137'
$ws.Range("E17").Value = '  +3.84%  '
$ws.Range("D18").Value = '2.621.97'
$ws.Range("E18").Value = '  +2.09%  '
$ws.Range("D19").Value = '4.51'
$ws.Range("E19").Value = '  +1.75%  '
$ws.Range("D20").Value = '338.63'
$ws.Range("E20").Value = '  +1.71%  '
$ws.Range("D21").Value = '10.25'
$ws.Range("E21").Value = '  +2.85%  '
$ws.Range("D22").Value = '6.28'
$ws.Range("E22").Value = '  +2.86%  '
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").Value = '65.20'
$ws.Range("E24").Value = '  -1.83%  '
$ws.Range("D25").Value = '0.448'
$ws.Range("E25").Value = '  +7.52%  '
$ws.Range("D26").Value = '0.163'
$ws.Range("E26").Value = '  +3.44%  '
$ws.Range("D27").Value = '0.998'
$ws.Range("E27").Value = '  -0.43%  '
$ws.Range("D28").Value = '7.30'
$ws.Range("E28").Value = '  +4.58%  '
$ws.Range("D29").Value = '0.0₃0780'
$ws.Range("E29").Value = '  +7.56%  '
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("D31").Value = '1.70'
$ws.Range("E31").Value = '  +3.43%  '
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").Value = '159.79'
$ws.Range("E32").Value = '  +3.21%  '
$ws.Range("B33").Value = 'Aptos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D33").Value = '6.04'
$ws.Range("E33").Value = '  +2.51%  '
$ws.Range("D34").Value = '18.99'
$ws.Range("E34").Value = '  +0.90%  '
$ws.Range("D35").Value = '4.06'
$ws.Range("E35").Value = '  +5.15%  '
$ws.Range("D36").Value = '0.891'
$ws.Range("E36").Value = '  +9.98%  '
$ws.Range("D37").Value = '1.14'
$ws.Range("E37").Value = '  +5.56%  '
$ws.Range("D38").Value = '0.870'
$ws.Range("E38").Value = '  +3.57%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").Value = '37.23'
$ws.Range("E39").Value = '  +0.93%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '1.49'
$ws.Range("E40").Value = '  +6.32%  '
$ws.Range("D41").Value = '294.86'
$ws.Range("E41").Value = '  +6.11%  '
$ws.Range("D42").Value = '3.62'
$ws.Range("E42").Value = '  +1.89%  '
$ws.Range("D43").Value = '0.997'
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").Value = '0.0977'
$ws.Range("E44").Value = '  +4.25%  '
$ws.Range("D45").Value = '0.594'
$ws.Range("E45").Value = '  +1.56%  '
$ws.Range("D46").Value = '0.0537'
$ws.Range("E46").Value = '  +1.85%  '
$ws.Range("B47").Value = 'WhiteBITCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D47").Value = '10.67'
$ws.Range("E47").Value = '  +0.36%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '19.08'
$ws.Range("E48").Value = '  +4.16%  '
$ws.Range("D49").Value = '125.87'
$ws.Range("E49").Value = '  +15.91%  '
$ws.Range("D50").Value = '0.0233'
$ws.Range("E50").Value = '  +3.74%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '1.941.24'
$ws.Range("E51").Value = '  +2.11%  '
